# Updates the pl_mw.xlsx result sheet for "case with 380 kV done".
# Replaces the numeric results in columns B,C,D,F,G,H,I,K,O for rows 2-25
# (the data rows of the single sheet) with their recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (as in the original file) that receive new values, in order.
$cols = @(2, 3, 4, 6, 7, 8, 9, 11, 15)   # B, C, D, F, G, H, I, K, O

# New values, one row per data row (rows 2..25), one column per entry in $cols above.
$data = @(
    @(0.3503875436509247,0.07152305426855321,0.03798645198736494,0.4029373145803277,0.2542690132076402,0.4497223053409449,0.3319660832775631,0.3733781911895733,1.327146005277669),
    @(0.3060370394515246,0.06574793501576437,0.03346670667245633,0.4073078043649652,0.2591425515850325,0.4551901566433436,0.3381545387602465,0.3263780134491867,1.348737384640145),
    @(0.278699529734098,0.06218614167966052,0.03067700653034677,0.410308906088197,0.2623812624152819,0.4587657458926984,0.3421849746020698,0.2973877949437451,1.362967551492744),
    @(0.2675334650646732,0.0607308201097112,0.02953659512822071,0.4116117073564283,0.2637629048532197,0.4602777373999096,0.3438853779610529,0.2855417985370536,1.369010932155142),
    @(0.2656778156855921,0.06048893544948442,0.02934701630366021,0.4118328568562966,0.2639960575849827,0.4605321192958272,0.3441712284998495,0.2835728594920965,1.370029190831154),
    @(0.2785490436110649,0.06216653017699514,0.03066164096357227,0.4103261529488229,0.262399645470321,0.4587859147907665,0.3422076721827896,0.2972281650323509,1.363048065029808),
    @(0.3351180275786021,0.06953515404821076,0.03643110896115331,0.4043783193053478,0.2558982207380645,0.4515623204272856,0.3340519559692083,0.3572004500454113,1.334388660621002),
    @(0.4451742309730093,0.08385453622885564,0.047626721883546,0.3952364480654715,0.2451080160695547,0.4391282335984172,0.3198908106133267,0.4737239870600547,1.285915762620874),
    @(0.5254615914982423,0.09429007978312143,0.05577702495135384,0.390059948469478,0.2383807654404322,0.4310470944409666,0.3106055530234189,0.5586349274441318,1.255024269204135),
    @(0.561855001251331,0.09901796781205974,0.05946790539057645,0.3880400040609189,0.2355822715122322,0.4275993994273364,0.3066248232634754,0.5971035251382375,1.241998325336858),
    @(0.5756167958982985,0.1008054101561413,0.06086307222010134,0.3873233007568331,0.2345603038385136,0.4263266842350149,0.3051524377725805,0.6116470554299269,1.237213586703845),
    @(0.5726538339449121,0.1004205838182202,0.06056271012765535,0.3874755108627781,0.2347787215299064,0.4265993251120079,0.3054679833689518,0.6085159155797442,1.238237483646984),
    @(0.562987590603143,0.09916508076788944,0.05958273697551419,0.3879800740413728,0.2354974361442856,0.4274940339756235,0.306502986809619,0.5983005106778592,1.241601716332873),
    @(0.5570641569190116,0.09839566651464793,0.05898214844729921,0.388295412886805,0.2359425907552897,0.4280463473113585,0.3071415201205667,0.592040165116174,1.243681675020568),
    @(0.5230805053805057,0.09398070166675154,0.05553547340075227,0.3901986995072733,0.2385689312543207,0.4312770065697151,0.3108706017470881,0.5561176529832323,1.255896228250251),
    @(0.5021987382965563,0.09126722827906519,0.05341670486357941,0.3914521100177524,0.2402472361378791,0.4333174246953782,0.3132206007158818,0.5340392230971815,1.263652621119206),
    @(0.4901759453647685,0.0897047020501418,0.05219647519996329,0.3922045552566225,0.241237184311089,0.4345125225464059,0.3145951436196759,0.5213254993925034,1.26821052060248),
    @(0.486103175915531,0.08917535100418661,0.05178305928355087,0.3924647317031429,0.2415765905894105,0.4349208553447603,0.3150644689841418,0.517018340480945,1.269770335004644),
    @(0.504422902690834,0.09155627024249213,0.05364241448702955,0.391315420093953,0.2400660275706485,0.4330979931871681,0.31296807028908,0.5363910481026153,1.262816937402363),
    @(0.5658273400511291,0.09953393237157115,0.05987064707299794,0.3878305629082917,0.2352853064638012,0.42723034469509,0.3061980300056701,0.6013016736837926,1.240609543538731),
    @(0.6058440095919195,0.1047307950926353,0.06392660779955861,0.3858340105856257,0.2323809849145775,0.4235869844591065,0.3019776299686168,0.6435859589345512,1.226957941121739),
    @(0.5844971653326922,0.1019587306660128,0.06176322400698098,0.3868738795961946,0.2339108924875859,0.4255139931600951,0.304211431230879,0.6210310476998018,1.234165081135146),
    @(0.5034174126767823,0.091425602191606,0.05354037776784537,0.3913771184260142,0.2401478738880236,0.4331971295377492,0.3130821661087362,0.5353278519219202,1.263194442759477),
    @(0.4154986941949517,0.07999528855910398,0.04461098617517223,0.3974393084682362,0.2478166951839498,0.4423067639049272,0.3235254210735317,0.4423211902054618,1.29820054766072)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowIndex, $cols[$j]).Value = $rowValues[$j]
    }
}
